{"js": "// Update benchmark stats table: fix README.md stats / docx prep values.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Simple single-value cell updates (row index -> new text)\nconst simpleUpdates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"202\",\n  4: \"0.00002\",\n  5: \"0.00006\",\n  8: \"0.00004\",\n  9: \"0.00004\",\n  11: \"0.00668\",\n};\n\nfor (const rowIndex in simpleUpdates) {\n  table.getCell(Number(rowIndex), 0).value = simpleUpdates[rowIndex];\n}\n\n// Rows whose multi-run (tab-separated) contents collapse into one value\nconst collapsedUpdates = {\n  43: \"100\",\n  44: \"0.01\",\n  45: \"147\",\n};\n\nfor (const rowIndex in collapsedUpdates) {\n  table.getCell(Number(rowIndex), 0).value = collapsedUpdates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table: fix README.md stats / docx prep values.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Simple single-value cell updates (1-based row index -> new text)\n$simpleUpdates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"202\"\n    5  = \"0.00002\"\n    6  = \"0.00006\"\n    9  = \"0.00004\"\n    10 = \"0.00004\"\n    12 = \"0.00668\"\n}\n\nforeach ($row in $simpleUpdates.Keys) {\n    $t.Cell($row, 1).Range.Text = $simpleUpdates[$row]\n}\n\n# Rows whose multi-run (tab-separated) contents collapse into one value\n$collapsedUpdates = @{\n    44 = \"100\"\n    45 = \"0.01\"\n    46 = \"147\"\n}\n\nforeach ($row in $collapsedUpdates.Keys) {\n    $t.Cell($row, 1).Range.Text = $collapsedUpdates[$row]\n}\n"}
